$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6 values ---
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "ClassPlus"
$ws.Range("C6").Value = 1288731

# D6 needs the same date number format as the other "Ticket Date" cells (D2:D5).
# Copy the format from D2 first, then set the value so the style (s="6") matches.
$ws.Range("D2").Copy($ws.Range("D6"))
$ws.Range("D6").Value = 45211

$ws.Range("E6").Value = "Anirban Chakraborty"
$ws.Range("F6").Value = "Cannot Manage Attandence in Batch"
$ws.Range("G6").Value = "Open"
$ws.Range("I6").Value = "Screenshot Sent"

# --- Re-style C6's font: was Bold 8pt Segoe UI -> now 11pt Oxygen (not bold) ---
$ws.Range("C6").Font.Bold = $false
$ws.Range("C6").Font.Size = 11
$ws.Range("C6").Font.Name = "Oxygen"

# --- Update the active selection to I6 ---
$ws.Range("I6").Select()
